$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Row 10 (rule "R30"): column C holds the "From" value.
# Restore/change it from 18 to 1, per the target revision.
$ws.Range("C10").Value = 1
